$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.165.33"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "'2.323.03"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'303.72"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'97.63"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").Value = "'35.55"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").Value = "'19.53"
$ws.Range("E11").Value = "  +7.94%  "
$ws.Range("D12").Value = "'0.0798"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").Value = "'2.687.00"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "'2.326.72"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "'42.974.11"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'12.61"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").Value = "'6.06"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'67.85"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'236.63"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  +3.84%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.44"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'24.91"
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.07"
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'165.19"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'9.14"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'33.23"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'18.03"
$ws.Range("E33").Value = "  +5.89%  "
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").Value = "'4.53"
$ws.Range("E35").Value = "  -8.16%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0697"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'2.34"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("D40").Value = "'1.76"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "'0.109"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "'1.991.08"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").Value = "'10.71"
$ws.Range("E43").Value = "  +5.99%  "
$ws.Range("D44").Value = "'0.0280"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").Value = "'17.99"
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").Value = "'2.07"
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("D47").Value = "'2.77"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "'2.89"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'2.554.13"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "'53.74"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").Value = "'71.93"
$ws.Range("E51").Value = "  -0.62%  "
